# Insert a new row at row 187. This pushes the existing rows 187-237
# down to 188-238, preserving all of their data/formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(187).Insert()

# Populate the newly inserted row 187 with its new record.
$ws.Range("A187").Value = 10
$ws.Range("B187").Value = "Vega Modelo de Temuco"
$ws.Range("C187").Value = "La Araucanía"
$ws.Range("D187").Value = 45204
$ws.Range("E187").Value = 9
$ws.Range("F187").Value = 100114002
$ws.Range("G187").Value = "Camote"
$ws.Range("H187").Value = "Sin especificar"
$ws.Range("I187").Value = "Primera"
$ws.Range("J187").Value = 120
$ws.Range("K187").Value = 24000
$ws.Range("L187").Value = 24000
$ws.Range("M187").Value = 24000
$ws.Range("N187").Value = "$/caja 18 kilos"
$ws.Range("O187").Value = "Perú"
$ws.Range("P187").Value = 1333
$ws.Range("Q187").Value = 18
$ws.Range("R187").Value = "Hortaliza"

# Make sure the date column keeps its date number format like the rest
# of column D.
$ws.Range("D187").NumberFormat = $ws.Range("D188").NumberFormat
